$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: A2 date 45242 -> 45254 ---
$ws.Range("A2").Value = 45254

# --- Rename header E1 (was "RSPM") to "PM10"; D1 stays "PM2" ---
$ws.Range("E1").Value = "PM10"

# --- Remove now-unused columns F:H (CO, O3, NH3) entirely ---
$ws.Range("F1:H11").Delete()

# --- Header row formatting: bold, size 12, row height 15.75 ---
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15.75

# --- Selection moves to E1 ---
$ws.Range("E1").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
